# Apply F-column ("想去人数" / want-to-go count) updates across all 4 sheets
# per commit "Update gh-pages to output generated at 456a3b4"
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 8
$ws.Cells.Item(4, 6).Value = 142
$ws.Cells.Item(6, 6).Value = 395
$ws.Cells.Item(7, 6).Value = 4890
$ws.Cells.Item(8, 6).Value = 4890
$ws.Cells.Item(9, 6).Value = 37
$ws.Cells.Item(11, 6).Value = 467
$ws.Cells.Item(13, 6).Value = 1114
$ws.Cells.Item(14, 6).Value = 656
$ws.Cells.Item(15, 6).Value = 4587
$ws.Cells.Item(16, 6).Value = 185
$ws.Cells.Item(18, 6).Value = 85
$ws.Cells.Item(20, 6).Value = 3601
$ws.Cells.Item(21, 6).Value = 7
$ws.Cells.Item(24, 6).Value = 3339
$ws.Cells.Item(26, 6).Value = 141
$ws.Cells.Item(28, 6).Value = 347
$ws.Cells.Item(30, 6).Value = 214
$ws.Cells.Item(32, 6).Value = 95
$ws.Cells.Item(33, 6).Value = 77
$ws.Cells.Item(34, 6).Value = 31
$ws.Cells.Item(37, 6).Value = 5841
$ws.Cells.Item(38, 6).Value = 916
$ws.Cells.Item(43, 6).Value = 1194
$ws.Cells.Item(44, 6).Value = 541
$ws.Cells.Item(46, 6).Value = 2065
$ws.Cells.Item(48, 6).Value = 76
$ws.Cells.Item(49, 6).Value = 731
$ws.Cells.Item(50, 6).Value = 873

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(24, 6).Value = 764
$ws.Cells.Item(25, 6).Value = 6

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 208

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 208
$ws.Cells.Item(7, 6).Value = 142
$ws.Cells.Item(9, 6).Value = 395
$ws.Cells.Item(10, 6).Value = 4890
$ws.Cells.Item(11, 6).Value = 4890
$ws.Cells.Item(12, 6).Value = 37
$ws.Cells.Item(15, 6).Value = 467
$ws.Cells.Item(16, 6).Value = 1114
$ws.Cells.Item(17, 6).Value = 656
$ws.Cells.Item(18, 6).Value = 4587
$ws.Cells.Item(19, 6).Value = 185
$ws.Cells.Item(21, 6).Value = 85
$ws.Cells.Item(23, 6).Value = 3601
$ws.Cells.Item(24, 6).Value = 3339
$ws.Cells.Item(26, 6).Value = 141
$ws.Cells.Item(27, 6).Value = 214
$ws.Cells.Item(29, 6).Value = 95
$ws.Cells.Item(30, 6).Value = 77
$ws.Cells.Item(35, 6).Value = 5841
$ws.Cells.Item(36, 6).Value = 916
$ws.Cells.Item(43, 6).Value = 1194
$ws.Cells.Item(44, 6).Value = 541
$ws.Cells.Item(45, 6).Value = 2065
$ws.Cells.Item(47, 6).Value = 76
$ws.Cells.Item(48, 6).Value = 731
$ws.Cells.Item(49, 6).Value = 873
